# "break out stock.yaml completed"
# 1) On the "day" sheet, the bsecode column (D) for rows 366-380 was being
#    stored as inline text; re-save those 15 cells as true numbers.
# 2) On the "week" sheet, 13 new rows (145-157) of stock data were appended.

$wb = $excel.ActiveWorkbook

# --- 1. Fix D366:D380 on "day" sheet: text -> number -------------------
$dayWs = $wb.Worksheets.Item("day")

$bsecodes = @{
    366 = 533179
    367 = 532644
    368 = 532504
    369 = 500520
    370 = 503806
    371 = 539957
    372 = 533273
    373 = 500670
    374 = 509480
    375 = 532514
    376 = 500104
    377 = 535755
    378 = 500477
    379 = 532155
    380 = 532754
}

foreach ($r in $bsecodes.Keys) {
    $dayWs.Cells.Item($r, 4).Value = $bsecodes[$r]
}

# --- 2. Append new rows to "week" sheet ---------------------------------
$weekWs = $wb.Worksheets.Item("week")

$newRows = @(
    @(1,  "HAL",        "Hindustan Aeronautics Ltd",          "541154", 2.32, 4769.8,  2878549,  "week", "16/08/2024 11:33:24"),
    @(2,  "RELIANCE",   "Reliance Industries Limited",        "500325", 1.12, 2956.4,  4708452,  "week", "16/08/2024 11:33:24"),
    @(3,  "NESTLEIND",  "Nestle India Limited",               "500790", 2.05, 2525.45, 859837,   "week", "16/08/2024 11:33:24"),
    @(4,  "KOTAKBANK",  "Kotak Mahindra Bank Limited",        "500247", 1.68, 1777.3,  2961871,  "week", "16/08/2024 11:33:24"),
    @(5,  "ICICIBANK",  "Icici Bank Limited",                 "532174", 2.2,  1187.25, 10961370, "week", "16/08/2024 11:33:24"),
    @(6,  "AXISBANK",   "Axis Bank Limited",                  "532215", 1.19, 1166.85, 6107134,  "week", "16/08/2024 11:33:24"),
    @(7,  "SBICARD",    "SBI Cards & Payment Services Ltd",   "543066", 1.31, 698.65,  815795,   "week", "16/08/2024 11:33:24"),
    @(8,  "RECLTD",     "Rural Electrification Corporation Limited", "532955", 2.83, 579.65, 6950751, "week", "16/08/2024 11:33:24"),
    @(9,  "PFC",        "Power Finance Corporation Limited",  "532810", 4.04, 504.25,  10378455, "week", "16/08/2024 11:33:24"),
    @(10, "LAURUSLABS", "Laurus Labs Limited",                "540222", 0.23, 429.5,   3156658,  "week", "16/08/2024 11:33:24"),
    @(11, "BEL",        "Bharat Electronics Limited",         "500049", 3.27, 303.3,   21353408, "week", "16/08/2024 11:33:24"),
    @(12, "BHEL",       "Bharat Heavy Electricals Limited",   "500103", 2.26, 296.55,  10118811, "week", "16/08/2024 11:33:24"),
    @(13, "NATIONALUM", "National Aluminium Company Limited", "532234", 0.02, 165.13,  23614094, "week", "16/08/2024 11:33:24")
)

$startRow = 145
$i = 0
foreach ($row in $newRows) {
    $r = $startRow + $i
    $weekWs.Cells.Item($r, 1).Value = $row[0]
    $weekWs.Cells.Item($r, 2).Value = $row[1]
    $weekWs.Cells.Item($r, 3).Value = $row[2]
    $weekWs.Cells.Item($r, 4).Value = "'" + $row[3]
    $weekWs.Cells.Item($r, 4).Style = "Normal"
    $weekWs.Cells.Item($r, 5).Value = $row[4]
    $weekWs.Cells.Item($r, 6).Value = $row[5]
    $weekWs.Cells.Item($r, 7).Value = $row[6]
    $weekWs.Cells.Item($r, 8).Value = $row[7]
    $weekWs.Cells.Item($r, 9).Value = $row[8]
    $i = $i + 1
}
